$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the score (C) and reviews count (E) columns to remain plain text
# (they hold numeric-looking strings, e.g. "8.2" / "2,546", and must not be
# auto-converted to numbers).
$ws.Range("C2:C26").NumberFormat = "@"
$ws.Range("E2:E26").NumberFormat = "@"

# Row 2: Maison le Bac Paris Aparthotel
$ws.Range("A2").Value = "Maison le Bac Paris Aparthotel"
$ws.Range("B2").Value = "US`$3,735"
$ws.Range("C2").Value = "8.2"
$ws.Range("D2").Value = "Very Good"
$ws.Range("E2").Value = "612"

# Row 3: Austin's Saint Lazare Hotel
$ws.Range("A3").Value = "Austin's Saint Lazare Hotel"
$ws.Range("B3").Value = "US`$5,908"
$ws.Range("C3").Value = "8.1"
$ws.Range("D3").Value = "Very Good"
$ws.Range("E3").Value = "2,546"

# Row 4: Edgar Suites Montmartre - Paul Albert
$ws.Range("A4").Value = "Edgar Suites Montmartre - Paul Albert"
$ws.Range("B4").Value = "US`$4,822"
$ws.Range("C4").Value = "9.3"
$ws.Range("D4").Value = "Wonderful"
$ws.Range("E4").Value = "114"

# Row 5: Austin's Arts Et Metiers Hotel
$ws.Range("A5").Value = "Austin's Arts Et Metiers Hotel"
$ws.Range("B5").Value = "US`$6,626"
$ws.Range("C5").Value = "8.2"
$ws.Range("D5").Value = "Very Good"
$ws.Range("E5").Value = "2,032"

# Row 6: Enjoy Hostel
$ws.Range("A6").Value = "Enjoy Hostel"
$ws.Range("B6").Value = "US`$917"
$ws.Range("C6").Value = "6.9"
$ws.Range("D6").Value = "Review score"
$ws.Range("E6").Value = "5,685"

# Row 7: PARIS AUTHENTIC HOUSE, Entier 1920's villa métro Line 7
$ws.Range("A7").Value = "PARIS AUTHENTIC HOUSE, Entier 1920's villa métro Line 7"
$ws.Range("B7").Value = "US`$9,295"
$ws.Range("C7").Value = "9.4"
$ws.Range("D7").Value = "Wonderful"
$ws.Range("E7").Value = "10"

# Row 8: Villa Royale Montsouris
$ws.Range("A8").Value = "Villa Royale Montsouris"
$ws.Range("B8").Value = "US`$4,327"
$ws.Range("C8").Value = "7.6"
$ws.Range("D8").Value = "Good"
$ws.Range("E8").Value = "283"

# Row 9: City Inn Paris
$ws.Range("A9").Value = "City Inn Paris"
$ws.Range("B9").Value = "US`$835"
$ws.Range("C9").Value = "6.7"
$ws.Range("D9").Value = "Review score"
$ws.Range("E9").Value = "2,498"

# Row 10: Hotel 29 Lepic
$ws.Range("A10").Value = "Hotel 29 Lepic"
$ws.Range("B10").Value = "US`$4,759"
$ws.Range("C10").Value = "8.1"
$ws.Range("D10").Value = "Very Good"
$ws.Range("E10").Value = "1,995"

# Row 11: GuestReady - Charm and Confort in the 18th
$ws.Range("A11").Value = "GuestReady - Charm and Confort in the 18th"
$ws.Range("B11").Value = "US`$2,443"
$ws.Range("C11").ClearContents()
$ws.Range("D11").ClearContents()
$ws.Range("E11").ClearContents()

# Row 12: Appart'City Collection Paris Grande Bibliothèque
$ws.Range("A12").Value = "Appart'City Collection Paris Grande Bibliothèque"
$ws.Range("B12").Value = "US`$4,294"
$ws.Range("C12").Value = "7.3"
$ws.Range("D12").Value = "Good"
$ws.Range("E12").Value = "1,368"

# Row 13: Résidence des Poissonniers
$ws.Range("A13").Value = "Résidence des Poissonniers"
$ws.Range("B13").Value = "US`$3,388"
$ws.Range("C13").Value = "7.3"
$ws.Range("D13").Value = "Good"
$ws.Range("E13").Value = "33"

# Row 14: Hôtel Soft
$ws.Range("A14").Value = "Hôtel Soft"
$ws.Range("B14").Value = "US`$4,481"
$ws.Range("C14").Value = "7.7"
$ws.Range("D14").Value = "Good"
$ws.Range("E14").Value = "1,691"

# Row 15: Fauchon l'Hôtel Paris
$ws.Range("A15").Value = "Fauchon l'Hôtel Paris"
$ws.Range("B15").Value = "US`$19,379"
$ws.Range("C15").Value = "9.3"
$ws.Range("D15").Value = "Wonderful"
$ws.Range("E15").Value = "882"

# Row 16: La Maison Gobert Paris Hotel Particulier
$ws.Range("A16").Value = "La Maison Gobert Paris Hotel Particulier"
$ws.Range("B16").Value = "US`$6,809"
$ws.Range("C16").Value = "9.2"
$ws.Range("D16").Value = "Wonderful"
$ws.Range("E16").Value = "492"

# Row 17: Hôtel Crimée
$ws.Range("A17").Value = "Hôtel Crimée"
$ws.Range("B17").Value = "US`$2,998"
$ws.Range("C17").Value = "5.9"
$ws.Range("D17").Value = "Review score"
$ws.Range("E17").Value = "581"

# Row 18: PARIS AUTHENTIC HOUSE 9 minutes by METRO RER B Gentilly to Notre-Dame de Paris
$ws.Range("A18").Value = "PARIS AUTHENTIC HOUSE 9 minutes by METRO RER B Gentilly to Notre-Dame de Paris"
$ws.Range("B18").Value = "US`$4,084"
$ws.Range("C18").Value = "7.0"
$ws.Range("D18").Value = "Good"
$ws.Range("E18").Value = "22"

# Row 19: Hollyday Studio Paris Centre Montmarte Sacré-coeur Opera Louvre
$ws.Range("A19").Value = "Hollyday Studio Paris Centre Montmarte Sacré-coeur Opera Louvre"
$ws.Range("B19").Value = "US`$4,338"
$ws.Range("C19").Value = "8.3"
$ws.Range("D19").Value = "Very Good"
$ws.Range("E19").Value = "27"

# Row 20: Hotel Darcet
$ws.Range("A20").Value = "Hotel Darcet"
$ws.Range("B20").Value = "US`$3,908"
$ws.Range("C20").Value = "8.8"
$ws.Range("D20").Value = "Excellent"
$ws.Range("E20").Value = "2,052"

# Row 21: CMG - Grands boulevards / Rex
$ws.Range("A21").Value = "CMG - Grands boulevards / Rex"
$ws.Range("B21").Value = "US`$6,333"
$ws.Range("C21").Value = "6.8"
$ws.Range("D21").Value = "Review score"
$ws.Range("E21").Value = "39"

# Row 22: Appartement Place du Trocadéro
$ws.Range("A22").Value = "Appartement Place du Trocadéro"
$ws.Range("B22").Value = "US`$8,576"
$ws.Range("C22").Value = "9.0"
$ws.Range("D22").Value = "Wonderful"
$ws.Range("E22").Value = "68"

# Row 23: LUXURY FLAT WITH ROOFTOP TERRACE - Paris 18
$ws.Range("A23").Value = "LUXURY FLAT WITH ROOFTOP TERRACE - Paris 18"
$ws.Range("B23").Value = "US`$10,065"
$ws.Range("C23").Value = "8.2"
$ws.Range("D23").Value = "Very Good"
$ws.Range("E23").Value = "27"

# Row 24: Best Stay Jeuneurs
$ws.Range("A24").Value = "Best Stay Jeuneurs"
$ws.Range("B24").Value = "US`$17,598"
$ws.Range("C24").Value = "8.3"
$ws.Range("D24").Value = "Very Good"
$ws.Range("E24").Value = "538"

# Row 25: Rent a Room - Residence Caire, Montorgueil
$ws.Range("A25").Value = "Rent a Room - Residence Caire, Montorgueil"
$ws.Range("B25").Value = "US`$7,368"
$ws.Range("C25").Value = "7.6"
$ws.Range("D25").Value = "Good"
$ws.Range("E25").Value = "88"

# Row 26: Apartments FOCH CHAMPS ELYSEES PARIS
$ws.Range("A26").Value = "Apartments FOCH CHAMPS ELYSEES PARIS"
$ws.Range("B26").Value = "US`$40,355"
$ws.Range("C26").Value = "7.8"
$ws.Range("D26").Value = "Good"
$ws.Range("E26").Value = "4"

# Rows 27-28 from the old data set are gone; drop them so the sheet
# ends at row 26 (dimension A1:E26).
$ws.Rows("27:28").Delete()
